$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "民國    年    月" placeholder text that used to live in the
# merged cell E3:G3 (data-period header). This also makes BA1/BB1/BC1/BD1
# (which derive the ROC year/month from E3) recompute to blanks/#VALUE!.
$ws.Range("E3").Value = ""

# Fill in explicit zero values for the new-origination figures that were
# previously left blank (rows 7-13, columns C-J).
$ws.Range("C7:J13").Value = 0

# Update the view: zoom to 85% and move the selection to F12.
$excel.ActiveWindow.Zoom = 85
$ws.Range("F12").Select()
